$d = $word.ActiveDocument

# 1. Change the title paragraph's style from Heading1 to Title
$p1 = $d.Paragraphs(1)
$p1.Style = "Title"

# 2. Replace the body paragraph's text with "AFFIDAVIT"
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "AFFIDAVIT"

# 3. Insert the remaining paragraphs after paragraph 2, one at a time,
#    walking forward so each InsertParagraphAfter lands after the
#    paragraph just written.
$cur = $p2

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "I, [Name], aged [Age] years, son/daughter of [Parent’s Name], residing at [Address], do hereby solemnly affirm and state as follows:"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "1. That I am the deponent herein and competent to swear this affidavit."

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "2. That I am making this affidavit to declare [insert reason]."

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "3. That the statements made herein are true to my knowledge and belief."

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Verified at [City] on this [Date]."

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "Signature:"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($cur.Index + 1)
$cur.Range.Text = "[Name]"

Write-Host "Done. Paragraph count: " $d.Paragraphs.Count
